# This workbook is a weekly price log for "Brócoli" at Vega Central Mapocho
# de Santiago. Each week contributes 4 rows (Primera/RM, Primera/O'Higgins,
# Segunda/RM, Segunda/O'Higgins) inserted at the top of the data block
# (row 252), pushing all older rows down by 4 and growing the used range
# from A1:R366 to A1:R370.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows at the top of the data block (row 252), shifting the
# existing rows 252:366 down to 256:370.
$ws.Range("A252:R255").EntireRow.Insert()

# New week's data: Fecha serial 44460 (2021-09-21)

# Row 252: Brócoli, Primera, Región Metropolitana
$ws.Range("A252").Value = 9
$ws.Range("B252").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C252").Value = "Metropolitana"
$ws.Range("D252").Value = 44460
$ws.Range("E252").Value = 13
$ws.Range("F252").Value = 100112023
$ws.Range("G252").Value = "Brócoli"
$ws.Range("H252").Value = "Sin especificar"
$ws.Range("I252").Value = "Primera"
$ws.Range("J252").Value = 4300
$ws.Range("K252").Value = 600
$ws.Range("L252").Value = 650
$ws.Range("M252").Value = 625
$ws.Range("N252").Value = "`$/unidad"
$ws.Range("O252").Value = "Región Metropolitana"
$ws.Range("P252").Value = 625
$ws.Range("Q252").Value = 1
$ws.Range("R252").Value = "Hortaliza"

# Row 253: Brócoli, Primera, Región de O'Higgins
$ws.Range("A253").Value = 9
$ws.Range("B253").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C253").Value = "Metropolitana"
$ws.Range("D253").Value = 44460
$ws.Range("E253").Value = 13
$ws.Range("F253").Value = 100112023
$ws.Range("G253").Value = "Brócoli"
$ws.Range("H253").Value = "Sin especificar"
$ws.Range("I253").Value = "Primera"
$ws.Range("J253").Value = 3400
$ws.Range("K253").Value = 600
$ws.Range("L253").Value = 650
$ws.Range("M253").Value = 625
$ws.Range("N253").Value = "`$/unidad"
$ws.Range("O253").Value = "Región de O'Higgins"
$ws.Range("P253").Value = 625
$ws.Range("Q253").Value = 1
$ws.Range("R253").Value = "Hortaliza"

# Row 254: Brócoli, Segunda, Región Metropolitana
$ws.Range("A254").Value = 9
$ws.Range("B254").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C254").Value = "Metropolitana"
$ws.Range("D254").Value = 44460
$ws.Range("E254").Value = 13
$ws.Range("F254").Value = 100112023
$ws.Range("G254").Value = "Brócoli"
$ws.Range("H254").Value = "Sin especificar"
$ws.Range("I254").Value = "Segunda"
$ws.Range("J254").Value = 2140
$ws.Range("K254").Value = 450
$ws.Range("L254").Value = 500
$ws.Range("M254").Value = 475
$ws.Range("N254").Value = "`$/unidad"
$ws.Range("O254").Value = "Región Metropolitana"
$ws.Range("P254").Value = 475
$ws.Range("Q254").Value = 1
$ws.Range("R254").Value = "Hortaliza"

# Row 255: Brócoli, Segunda, Región de O'Higgins
$ws.Range("A255").Value = 9
$ws.Range("B255").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C255").Value = "Metropolitana"
$ws.Range("D255").Value = 44460
$ws.Range("E255").Value = 13
$ws.Range("F255").Value = 100112023
$ws.Range("G255").Value = "Brócoli"
$ws.Range("H255").Value = "Sin especificar"
$ws.Range("I255").Value = "Segunda"
$ws.Range("J255").Value = 1330
$ws.Range("K255").Value = 450
$ws.Range("L255").Value = 500
$ws.Range("M255").Value = 475
$ws.Range("N255").Value = "`$/unidad"
$ws.Range("O255").Value = "Región de O'Higgins"
$ws.Range("P255").Value = 475
$ws.Range("Q255").Value = 1
$ws.Range("R255").Value = "Hortaliza"
